$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue $ws 'D2' '65.909.45'
$ws.Range('E2').Value = '  -5.38%  '

# Row 3
Set-TextValue $ws 'D3' '3.323.09'
$ws.Range('E3').Value = '  -6.42%  '

# Row 4
$ws.Range('E4').Value = '  +0.20%  '

# Row 5
Set-TextValue $ws 'D5' '558.44'
$ws.Range('E5').Value = '  -4.67%  '

# Row 6
Set-TextValue $ws 'D6' '180.68'
$ws.Range('E6').Value = '  -8.30%  '

# Row 7
$ws.Range('E7').Value = '  +0.05%  '

# Row 8
Set-TextValue $ws 'D8' '0.587'
$ws.Range('E8').Value = '  -4.02%  '

# Row 9
Set-TextValue $ws 'D9' '3.311.90'
$ws.Range('E9').Value = '  -6.38%  '

# Row 10
Set-TextValue $ws 'D10' '0.183'
$ws.Range('E10').Value = '  -11.43%  '

# Row 11
Set-TextValue $ws 'D11' '0.582'
$ws.Range('E11').Value = '  -7.80%  '

# Row 12
Set-TextValue $ws 'D12' '47.08'
$ws.Range('E12').Value = '  -9.65%  '

# Row 13
Set-TextValue $ws 'D13' '0.0000262'
$ws.Range('E13').Value = '  -8.95%  '

# Row 14
Set-TextValue $ws 'D14' '3.861.33'

# Row 15
$ws.Range('E15').Value = '  -8.07%  '

# Row 16
Set-TextValue $ws 'D16' '598.88'
$ws.Range('E16').Value = '  -10.77%  '

# Row 17
Set-TextValue $ws 'D17' '65.904.32'
$ws.Range('E17').Value = '  -5.50%  '

# Row 18
Set-TextValue $ws 'D18' '17.95'
$ws.Range('E18').Value = '  -2.81%  '

# Row 19
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws 'D19' '0.117'
$ws.Range('E19').Value = '  -3.96%  '

# Row 20
$ws.Range('B20').Value = 'WrappedEther'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws 'D20' '3.318.03'
$ws.Range('E20').Value = '  -6.53%  '

# Row 21
Set-TextValue $ws 'D21' '11.35'
$ws.Range('E21').Value = '  -9.36%  '

# Row 22
Set-TextValue $ws 'D22' '0.900'
$ws.Range('E22').Value = '  -6.83%  '

# Row 23
Set-TextValue $ws 'D23' '16.72'
$ws.Range('E23').Value = '  -7.31%  '

# Row 24
Set-TextValue $ws 'D24' '5.03'
$ws.Range('E24').Value = '  -6.16%  '

# Row 25
Set-TextValue $ws 'D25' '99.51'
$ws.Range('E25').Value = '  -5.73%  '

# Row 26
Set-TextValue $ws 'D26' '3.99'
$ws.Range('E26').Value = '  -9.13%  '

# Row 27
Set-TextValue $ws 'D27' '6.00'
$ws.Range('E27').Value = '  -0.19%  '

# Row 28
Set-TextValue $ws 'D28' '2.63'
$ws.Range('E28').Value = '  -10.07%  '

# Row 29
Set-TextValue $ws 'D29' '9.21'
$ws.Range('E29').Value = '  -9.54%  '

# Row 30
Set-TextValue $ws 'D30' '8.61'
$ws.Range('E30').Value = '  -10.51%  '

# Row 31
$ws.Range('E31').Value = '  -9.04%  '

# Row 32
$ws.Range('E32').Value = '  -8.41%  '

# Row 33
Set-TextValue $ws 'D33' '3.76'
$ws.Range('E33').Value = '  -14.09%  '

# Row 34
Set-TextValue $ws 'D34' '10.95'
$ws.Range('E34').Value = '  -6.83%  '

# Row 35
$ws.Range('E35').Value = '  -6.89%  '

# Row 36
Set-TextValue $ws 'D36' '3.752.96'
$ws.Range('E36').Value = '  -0.78%  '

# Row 37
Set-TextValue $ws 'D37' '57.69'
$ws.Range('E37').Value = '  -7.04%  '

# Row 38
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws 'D38' '528.82'
$ws.Range('E38').Value = '  +5.30%  '

# Row 39
$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws 'D39' '0.999'
$ws.Range('E39').Value = '  -0.18%  '

# Row 40
Set-TextValue $ws 'D40' '3.44'
$ws.Range('E40').Value = '  -7.76%  '

# Row 41
Set-TextValue $ws 'D41' '0.0₃0708'
$ws.Range('E41').Value = '  -13.49%  '

# Row 42
$ws.Range('E42').Value = '  -9.07%  '

# Row 43
$ws.Range('E43').Value = '  -8.08%  '

# Row 44
Set-TextValue $ws 'D44' '0.337'
$ws.Range('E44').Value = '  -9.22%  '

# Row 45
Set-TextValue $ws 'D45' '31.69'
$ws.Range('E45').Value = '  -8.61%  '

# Row 46
Set-TextValue $ws 'D46' '3.26'
$ws.Range('E46').Value = '  -3.62%  '

# Row 47
Set-TextValue $ws 'D47' '0.0411'
$ws.Range('E47').Value = '  -8.64%  '

# Row 48
$ws.Range('B48').Value = 'CoreDAO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
Set-TextValue $ws 'D48' '2.98'
$ws.Range('E48').Value = '  +9.20%  '

# Row 49
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws 'D49' '0.129'
$ws.Range('E49').Value = '  -5.90%  '

# Row 50
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue $ws 'D50' '2.59'
$ws.Range('E50').Value = '  -9.86%  '

# Row 51
Set-TextValue $ws 'D51' '0.999'
$ws.Range('E51').Value = '  -0.04%  '
